$wb = $excel.ActiveWorkbook

# Source sheet used as the template for the three new market sheets -
# it already carries the right column widths / row heights / styles.
$src = $wb.Worksheets.Item("UK")

# --- Denmark ---------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$wsDenmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDenmark.Name = "Denmark"
$wsDenmark.Range("B2").Value = "Denmark Market"
$wsDenmark.Range("B4").Value = "NGC-3446/T2008"
$wsDenmark.Range("A1:XFD1048576").Select()

# --- Sweden ------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$wsSweden = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSweden.Name = "Sweden"
$wsSweden.Range("B2").Value = "Sweden Market"
$wsSweden.Range("B4").Value = "NGC-3465/T2023"
$wsSweden.Range("A1:XFD1048576").Select()

# --- Norway --------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$wsNorway = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNorway.Name = "Norway"
$wsNorway.Range("B2").Value = "Norway Market"
$wsNorway.Range("B4").Value = "NGC-3464/T1920"
$wsNorway.Range("A9").Select()

$wsNorway.Activate()
